$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Underline the existing section headers (CARBON ONLY / CARBON and SILICON)
# to match the new "C, Si, N, Al" section header style.
$ws.Range("A2").Font.Underline = $true
$ws.Range("A18").Font.Underline = $true

# New section: "C, Si, N, Al" isotopes model ratings.
$ws.Range("A26").Value = "C, Si, N, Al"
$ws.Range("A26").Font.Underline = $true

$ws.Range("A27").Value = "kNN"
$ws.Range("B27").Value = 85.4

$ws.Range("A28").Value = "Logistic Regression"
$ws.Range("B28").Value = 88.5

$ws.Range("A29").Value = "Decision Tree"
$ws.Range("B29").Value = 93.9

$ws.Range("A30").Value = "Random Forest Classifier"
$ws.Range("B30").Value = 95.1

# Scroll the view down to the new section and restore the original selection.
$ws.Range("C29").Select()
